$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: merge R27 into the existing R6, R28 resistor line ---
# Qty 2 -> 3
$ws.Range("A3").Value = 3
# Parts "R6, R28" -> "R6, R27, R28"
$ws.Range("E3").Value = "R6, R27, R28"

# --- New row 32: SW3 slide toggle switch (DPDT) ---
$ws.Range("A32").Value = 1
$ws.Range("B32").Value = "JS202011JCQN"
$ws.Range("C32").Value = "JS202011JCQN"
$ws.Range("D32").Value = "JS202011JCQN"
$ws.Range("E32").Value = "SW3"
$ws.Range("F32").Value = "SLIDE TOGGLE SWITCH VERTICAL (DPDT)"
$ws.Range("G32").Value = "C&K Components"
$ws.Range("H32").Value = "JS202011JCQN"
$ws.Range("I32").Value = "CKN10723CT-ND"
$ws.Range("J32").Value = "611-JS202011JCQN"

# Match the row's formatting (fontId 1 / explicit black rgb) used by the rest
# of the table's styled rows.
$ws.Range("A32:K32").Font.Color = 0

# --- Selection / view bookkeeping to match the final saved state ---
$ws.Range("D39").Select()
